$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-0.12, 0.147, -0.201, 91, 28, -160),
    @(0.008999999999999999, -0.295, 0.461, 5, -54, 49),
    @(-0.191, -0.205, 0.617, -157, 0, -114),
    @(-0.154, 0.259, 0.615, 86, 39, -75),
    @(0.165, -0.059, 0.846, -3, -11, -92),
    @(0.108, -0.165, 0.5679999999999999, 113, -3, 10),
    @(-0.492, 0.095, 0.36, -89, -19, 149),
    @(-0.1, 0.485, 0.734, 86, 7, -152),
    @(0.046, -0.309, 0.168, -32, 41, -71),
    @(-0.278, -0.252, 0.634, 126, 54, -19)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
